$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1124.2059
$ws.Range("I28").Value = 935.13336
$ws.Range("J28").Value = 2542.25
$ws.Range("K28").Value = 935.13336
$ws.Range("L28").Value = 2542.25
$ws.Range("M28").Value = -450.13336
$ws.Range("N28").Value = -3512.25

$ws.Range("H64").Value = 66669656
$ws.Range("I64").Value = 66669656
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 66669656
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = $null
$ws.Range("N64").Value = -66669408

$ws.Range("H67").Value = 66669656
$ws.Range("I67").Value = 66669656
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 66669656
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = $null
$ws.Range("N67").Value = -66668798

$ws.Range("H76").Value = 7831
$ws.Range("I76").Value = 7249.5
$ws.Range("J76").Value = 8412.5
$ws.Range("K76").Value = 7249.5
$ws.Range("L76").Value = 8412.5
$ws.Range("M76").Value = -6934.5
$ws.Range("N76").Value = -9042.5

$ws.Range("H79").Value = 7831
$ws.Range("I79").Value = 7249.5
$ws.Range("J79").Value = 8412.5
$ws.Range("K79").Value = 7249.5
$ws.Range("L79").Value = 8412.5
$ws.Range("M79").Value = -6157.5
$ws.Range("N79").Value = -10596.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 100000000
$ws.Range("I102").Value = 100000000
$ws.Range("K102").Value = 100000000
$ws.Range("M102").Value = -99998378

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 100000
$ws.Range("J87").Value = 100000
$ws.Range("L87").Value = 100000
$ws.Range("N87").Value = -102496

$ws.Range("H90").Value = 100000
$ws.Range("J90").Value = 100000
$ws.Range("L90").Value = 300000
$ws.Range("N90").Value = -312480

$ws.Range("H99").Value = 1732.5834
$ws.Range("I99").Value = 900.3333
$ws.Range("J99").Value = 2010
$ws.Range("K99").Value = 900.3333
$ws.Range("L99").Value = 2010
$ws.Range("M99").Value = 597.6667
$ws.Range("N99").Value = -5006

$ws.Range("H107").Value = 72522.42999999999
$ws.Range("I107").Value = 1109.5834
$ws.Range("J107").Value = 500999.5
$ws.Range("K107").Value = 1109.5834
$ws.Range("L107").Value = 500999.5
$ws.Range("M107").Value = 810.4166
$ws.Range("N107").Value = -504839.5

$ws.Range("H134").Value = 29414794
$ws.Range("I134").Value = 31251968
$ws.Range("K134").Value = 93755904
$ws.Range("M134").Value = -93753369

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9363.755999999999
$ws.Range("J31").Value = 12021.952
$ws.Range("L31").Value = 12021.952
$ws.Range("N31").Value = -12611.952

$ws.Range("H34").Value = 9363.755999999999
$ws.Range("J34").Value = 12021.952
$ws.Range("L34").Value = 12021.952
$ws.Range("N34").Value = -12425.952

$ws.Range("H62").Value = 4939.8
$ws.Range("I62").Value = 4949.75
$ws.Range("J62").Value = 4900
$ws.Range("K62").Value = 4949.75
$ws.Range("L62").Value = 4900
$ws.Range("M62").Value = -4325.75
$ws.Range("N62").Value = -6148

$ws.Range("H65").Value = 4939.8
$ws.Range("I65").Value = 4949.75
$ws.Range("J65").Value = 4900
$ws.Range("K65").Value = 24748.75
$ws.Range("L65").Value = 24500
$ws.Range("M65").Value = -21628.75
$ws.Range("N65").Value = -30740

$ws.Range("H107").Value = 270780.78
$ws.Range("I107").Value = 388504.7
$ws.Range("K107").Value = 388504.7
$ws.Range("M107").Value = -386584.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = $null

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").Value = $null

$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 852748.7
$ws.Range("I7").Value = 1684664.9
$ws.Range("J7").Value = 20832.5
$ws.Range("K7").Value = 1684664.9
$ws.Range("L7").Value = 20832.5
$ws.Range("M7").Value = -1684552.9
$ws.Range("N7").Value = -21056.5

$ws.Range("H8").Value = 852748.7
$ws.Range("I8").Value = 1684664.9
$ws.Range("J8").Value = 20832.5
$ws.Range("K8").Value = 1684664.9
$ws.Range("L8").Value = 20832.5
$ws.Range("M8").Value = -1684525.9
$ws.Range("N8").Value = -21110.5

$ws.Range("H10").Value = 11200.6
$ws.Range("I10").Value = 13000.75
$ws.Range("J10").Value = 4000
$ws.Range("K10").Value = 13000.75
$ws.Range("L10").Value = 4000
$ws.Range("M10").Value = -12831.75
$ws.Range("N10").Value = -4338

$ws.Range("H11").Value = 6909.769
$ws.Range("I11").Value = 8899.700000000001
$ws.Range("K11").Value = 8899.700000000001
$ws.Range("M11").Value = -8760.700000000001

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = $null

$ws.Range("H14").Value = 2545002
$ws.Range("I14").Value = 9999999
$ws.Range("J14").Value = 60003
$ws.Range("K14").Value = 9999999
$ws.Range("L14").Value = 60003
$ws.Range("M14").Value = -9999831
$ws.Range("N14").Value = -60339

$ws.Range("H80").Value = 2216.1667
$ws.Range("I80").Value = 2099
$ws.Range("J80").Value = 2333.3333
$ws.Range("K80").Value = 2099
$ws.Range("L80").Value = 2333.3333
$ws.Range("M80").Value = -1101
$ws.Range("N80").Value = -4329.3333

$ws.Range("H83").Value = 2216.1667
$ws.Range("I83").Value = 2099
$ws.Range("J83").Value = 2333.3333
$ws.Range("K83").Value = 10495
$ws.Range("L83").Value = 11666.6665
$ws.Range("M83").Value = -5503
$ws.Range("N83").Value = -21650.6665

$ws.Range("H98").Value = 43650
$ws.Range("J98").Value = 43650
$ws.Range("L98").Value = 43650
$ws.Range("N98").Value = -49640

$ws.Range("H107").Value = 765.9091
$ws.Range("I107").Value = 333.57144
$ws.Range("K107").Value = 333.57144
$ws.Range("M107").Value = 1586.42856

$ws.Range("H113").Value = 61218.47
$ws.Range("I113").Value = 73360.57000000001
$ws.Range("J113").Value = 4555.3335
$ws.Range("K113").Value = 73360.57000000001
$ws.Range("L113").Value = 4555.3335
$ws.Range("M113").Value = -71190.57000000001
$ws.Range("N113").Value = -8895.333500000001

$ws.Range("H132").Value = 7357500.5
$ws.Range("I132").Value = 8336293.5
$ws.Range("J132").Value = 16555.5
$ws.Range("K132").Value = 25008880.5
$ws.Range("L132").Value = 49666.5
$ws.Range("M132").Value = -25006350.5
$ws.Range("N132").Value = -54726.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2263.5264
$ws.Range("I22").Value = 2148.8
$ws.Range("J22").Value = 2391
$ws.Range("K22").Value = 2148.8
$ws.Range("L22").Value = 2391
$ws.Range("M22").Value = -1853.8
$ws.Range("N22").Value = -2981

$ws.Range("H27").Value = 2263.5264
$ws.Range("I27").Value = 2148.8
$ws.Range("J27").Value = 2391
$ws.Range("K27").Value = 2148.8
$ws.Range("L27").Value = 2391
$ws.Range("M27").Value = -2041.8
$ws.Range("N27").Value = -2605

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 99999
$ws.Range("J130").Value = 99999
$ws.Range("L130").Value = 99999
$ws.Range("N130").Value = -110039
